$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (16.42578125 -> 15.42578125, i.e. one character narrower)
$ws.Columns.Item(1).ColumnWidth = 14.625

# New values for columns A and B, rows 1-32
$colA = @(-0.11035056273264843, -0.10139496781091228, -0.051533041652325196, -0.043376120604166957, -0.039963717565440149, -0.030260325930264287, -0.020015671884376474, -0.0099684123757968202, -0.007914285857660186, -0.027397979469913736, -0.024382843095390072, -0.020857179546180937, -0.017167009992888005, -0.0090804452707251571, -0.0080522136202301198, -0.0060339888265481534, -0.0040034072513561725, -0.0028982066194487288, 0.001170792801813203, 0.0056699560893562762, 0.0098070066763202846, -0.045701268415552221, -0.040491227942330177, -0.020097449563091452, -0.0063029559696730075, -0.0037473661263405234, -0.0011784509877026395, 0.0012137055584071632, 0.0077425992366801921, -0.021164315503148678, -0.014022348281992336, -0.004001218526784811)
$colB = @(0.11029317439567166, 0.10123431003982386, 0.051376120551628546, 0.042963717539567625, 0.038563281875277688, 0.030015671813698575, 0.019968412303628558, 0.0099142858296334957, 0.0078787223074940904, 0.027382843060440365, 0.024357179507537907, 0.020667009950566584, 0.017080445201885119, 0.0090522135904889112, 0.0080339887905322982, 0.0060034072144383721, 0.003999999951511235, 0.0028292071747486602, -0.0016699561129631668, -0.0058070066997650827, -0.010013515737603562, 0.045491227906605758, 0.04009744943842275, 0.019999999873585139, 0.0062473660987034663, 0.0036784509599474546, 0.00078629441637056985, -0.0014653052578417203, -0.0078120929289555718, 0.02102234822644844, 0.01400121845380653, 0.0039999999609374726)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}
